$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("receipt_multiple")

# Update the total_amount value on row 2 (G2): 1209 -> 0
$ws.Range("G2").Value = 0

# Update the ri_number value on row 6 (E6): 20999999999 -> 99999999999
$ws.Range("E6").Value = 99999999999

# Update the active selection cell from E7 to F7
$ws.Range("F7").Select()
